# New crime data collected - update 43rd Precinct weekly CompStat report
# (Volume/Number header, reporting week dates, and the weekly crime-complaint
# statistics table for rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Volume 31   Number  17" -> "...18" ---
$ws.Range("A8").Value = "Volume 31   Number  18"

# --- Header text: reporting week dates ---
$ws.Range("C9").Value = "Report Covering the Week  4/29/2024  Through  5/5/2024"

# --- Cells that switch from the blank-placeholder shared string ("0") to a
#     real numeric value now need the numeric (#,##0) format before the
#     value is written so the cell style matches a genuine number cell. ---
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("F22").NumberFormat = "#,##0"

# Row 14
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = -33.333333333333
$ws.Range("L14").Value = -50
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -92
# Row 15
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -75
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 19
$ws.Range("K15").Value = -47.368421052631
$ws.Range("L15").Value = -54.545454545454
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -60
# Row 16
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -27.272727272727
$ws.Range("F16").Value = 43
$ws.Range("G16").Value = 42
$ws.Range("H16").Value = 2.380952380952
$ws.Range("I16").Value = 164
$ws.Range("J16").Value = 184
$ws.Range("K16").Value = -10.869565217391
$ws.Range("L16").Value = -20
$ws.Range("M16").Value = 4.458598726114
$ws.Range("N16").Value = -75.917767988252
# Row 17
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -15.384615384615
$ws.Range("F17").Value = 54
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = -21.739130434782
$ws.Range("I17").Value = 222
$ws.Range("J17").Value = 275
$ws.Range("K17").Value = -19.272727272727
$ws.Range("L17").Value = 0.452488687782
$ws.Range("M17").Value = 30.588235294117
$ws.Range("N17").Value = -21.554770318021
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 18
$ws.Range("H18").Value = -28
$ws.Range("I18").Value = 98
$ws.Range("J18").Value = 130
$ws.Range("K18").Value = -24.615384615384
$ws.Range("L18").Value = -1.010101010101
$ws.Range("M18").Value = -22.834645669291
$ws.Range("N18").Value = -83.501683501683
# Row 19
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 96
$ws.Range("H19").Value = -29.166666666666
$ws.Range("I19").Value = 302
$ws.Range("J19").Value = 348
$ws.Range("K19").Value = -13.218390804597
$ws.Range("L19").Value = -14.204545454545
$ws.Range("M19").Value = 72.571428571428
$ws.Range("N19").Value = 16.153846153846
# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 15
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 62
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 152
$ws.Range("J20").Value = 230
$ws.Range("K20").Value = -33.913043478260
$ws.Range("L20").Value = -22.051282051282
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -76.969696969697
# Row 21
$ws.Range("C21").Value = 55
$ws.Range("D21").Value = 71
$ws.Range("E21").Value = -22.535211267605
$ws.Range("F21").Value = 217
$ws.Range("G21").Value = 302
$ws.Range("H21").Value = -28.145695364238
$ws.Range("I21").Value = 950
$ws.Range("J21").Value = 1189
$ws.Range("K21").Value = -20.100925147182
$ws.Range("L21").Value = -13.479052823315
$ws.Range("M21").Value = 32.496513249651
$ws.Range("N21").Value = -62.420886075949
# Row 22
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 16.666666666666
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = -12.5
# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -77.777777777777
$ws.Range("F23").Value = 20
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = -4.761904761904
$ws.Range("I23").Value = 86
$ws.Range("J23").Value = 101
$ws.Range("K23").Value = -14.851485148514
$ws.Range("L23").Value = -22.522522522522
$ws.Range("M23").Value = 24.637681159420
# Row 24
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -14.393939393939
$ws.Range("I24").Value = 582
$ws.Range("J24").Value = 673
$ws.Range("K24").Value = -13.521545319465
$ws.Range("L24").Value = -7.765451664025
$ws.Range("M24").Value = 23.829787234042
# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -63.636363636363
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 53
$ws.Range("H25").Value = -45.283018867924
$ws.Range("I25").Value = 199
$ws.Range("J25").Value = 290
$ws.Range("K25").Value = -31.379310344827
$ws.Range("L25").Value = -47.354497354497
# Row 26
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 37
$ws.Range("E26").Value = -32.432432432432
$ws.Range("F26").Value = 84
$ws.Range("G26").Value = 100
$ws.Range("H26").Value = -16
$ws.Range("I26").Value = 395
$ws.Range("J26").Value = 373
$ws.Range("K26").Value = 5.898123324396
$ws.Range("L26").Value = 5.614973262032
$ws.Range("M26").Value = -23.151750972762
# Row 27
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -60
$ws.Range("I27").Value = 18
$ws.Range("J27").Value = 26
$ws.Range("K27").Value = -30.769230769230
$ws.Range("L27").Value = -41.935483870967
# Row 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 200
$ws.Range("F28").Value = 13
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 225
$ws.Range("I28").Value = 43
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = 19.444444444444
$ws.Range("L28").Value = 104.761904761905
# Row 29
$ws.Range("G29").Value = 4
$ws.Range("L29").Value = -35.714285714285
$ws.Range("N29").Value = -80.434782608695
# Row 30
$ws.Range("G30").Value = 1
$ws.Range("L30").Value = -42.857142857142
$ws.Range("N30").Value = -81.395348837209
